$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "satara"
$ws.Range("C1").Value = "reshma"
$ws.Range("A2").Value = "gangapur"
$ws.Range("B2").Value = "didi"
$ws.Range("C2").Value = "kalamboli"
